# Updates Price (D) and Volume(1h) (E) columns for the cryptos table;
# rows 40/41 also get Coin (B) and Link (C) swapped per upstream reorder.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{Row=2; B='Bitcoin'; C='https://coinranking.com/coin/Qwsogvtv82FCd+bitcoin-btc'; D='28.320.73'; E='  +5.94%  '},
    @{Row=3; B='Ethereum'; C='https://coinranking.com/coin/razxDUgYGNAdQ+ethereum-eth'; D='1.790.06'; E='  +3.39%  '},
    @{Row=4; B='TetherUSD'; C='https://coinranking.com/coin/HIVsRcGKkPFtW+tetherusd-usdt'; D='0.9988'; E='  +0.09%  '},
    @{Row=5; B='BNB'; C='https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb'; D='247.29'; E='  +2.11%  '},
    @{Row=6; B='USDC'; C='https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc'; D='0.9985'; E='  +0.01%  '},
    @{Row=7; B='XRP'; C='https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp'; D='0.4925'; E='  +0.12%  '},
    @{Row=8; B='Cardano'; C='https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'; D='0.2695'; E='  +2.78%  '},
    @{Row=9; B='Dogecoin'; C='https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'; D='0.06294'; E='  +1.23%  '},
    @{Row=10; B='WrappedEther'; C='https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'; D='1.785.22'; E='  +3.10%  '},
    @{Row=11; B='Solana'; C='https://coinranking.com/coin/zNZHO_Sjf+solana-sol'; D='16.61'; E='  +4.08%  '},
    @{Row=12; B='TRON'; C='https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'; D='0.07053'; E='  +1.05%  '},
    @{Row=13; B='Polygon'; C='https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'; D='0.6305'; E='  +2.98%  '},
    @{Row=14; B='Polkadot'; C='https://coinranking.com/coin/25W7FG7om+polkadot-dot'; D='4.670'; E='  +3.53%  '},
    @{Row=15; B='Litecoin'; C='https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'; D='80.26'; E='  +3.93%  '},
    @{Row=16; B='WrappedBTC'; C='https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'; D='28.269.38'; E='  +6.56%  '},
    @{Row=17; B='Dai'; C='https://coinranking.com/coin/MoTuySvg7+dai-dai'; D='0.9992'; E='  +0.08%  '},
    @{Row=18; B='BinanceUSD'; C='https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'; D='0.9984'; E='  +0.03%  '},
    @{Row=19; B='ShibaInu'; C='https://coinranking.com/coin/xz24e0BjL+shibainu-shib'; D='0.000007265'; E='  +0.89%  '},
    @{Row=20; B='Avalanche'; C='https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'; D='12.08'; E='  +5.85%  '},
    @{Row=21; B='WrappedliquidstakedEther2.0'; C='https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'; D='2.015.93'; E='  +3.34%  '},
    @{Row=22; B='Uniswap'; C='https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'; D='4.569'; E='  +2.35%  '},
    @{Row=23; B='Cosmos'; C='https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'; D='8.781'; E='  +2.48%  '},
    @{Row=24; B='Chainlink'; C='https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'; D='5.262'; E='  +3.00%  '},
    @{Row=25; B='Monero'; C='https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'; D='141.95'; E='  +2.82%  '},
    @{Row=26; B='EthereumClassic'; C='https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'; D='15.79'; E='  +3.08%  '},
    @{Row=27; B='LidoDAOToken'; C='https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'; D='1.857'; E='  +5.05%  '},
    @{Row=28; B='BitcoinCash'; C='https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'; D='109.94'; E='  +3.30%  '},
    @{Row=29; B='Toncoin'; C='https://coinranking.com/coin/67YlI0K1b+toncoin-ton'; D='1.387'; E='  +0.23%  '},
    @{Row=30; B='InternetComputer(DFINITY)'; C='https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'; D='4.178'; E='  +6.09%  '},
    @{Row=31; B='Stellar'; C='https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'; D='0.08282'; E='  +3.79%  '},
    @{Row=32; B='Filecoin'; C='https://coinranking.com/coin/ymQub4fuB+filecoin-fil'; D='3.784'; E='  +3.07%  '},
    @{Row=33; B='Hedera'; C='https://coinranking.com/coin/jad286TjB+hedera-hbar'; D='0.04910'; E='  +9.48%  '},
    @{Row=34; B='ARBITRUM'; C='https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'; D='1.086'; E='  +8.26%  '},
    @{Row=35; B='ImmutableX'; C='https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'; D='0.6612'; E='  +6.16%  '},
    @{Row=36; B='HuobiToken'; C='https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'; D='2.614'; E='  +0.25%  '},
    @{Row=37; B='TrustWalletToken'; C='https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'; D='0.9479'; E='  +0.45%  '},
    @{Row=38; B='MXToken'; C='https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'; D='2.612'; E='  +7.71%  '},
    @{Row=39; B='RenderToken'; C='https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'; D='2.080'; E='  +1.40%  '},
    @{Row=40; B='VeChain'; C='https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'; D='0.01558'; E='  +3.24%  '},
    @{Row=41; B='FraxShare'; C='https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'; D='5.909'; E='  +6.06%  '},
    @{Row=42; B='PaxDollar'; C='https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'; D='0.9992'; E='  +0.15%  '},
    @{Row=43; B='Quant'; C='https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'; D='99.90'; E='  +0.37%  '},
    @{Row=44; B='TheSandbox'; C='https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'; D='0.4009'; E='  +3.85%  '},
    @{Row=45; B='Aptos'; C='https://coinranking.com/coin/HGYj5JCv5+aptos-apt'; D='7.213'; E='  +4.45%  '},
    @{Row=46; B='Algorand'; C='https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'; D='0.1214'; E='  +4.56%  '},
    @{Row=47; B='Cronos'; C='https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'; D='0.05453'; E='  +1.35%  '},
    @{Row=48; B='EnergySwap'; C='https://coinranking.com/coin/SbWqqTui-+energyswap-ens'; D='8.053'; E='  +2.05%  '},
    @{Row=49; B='Elrond'; C='https://coinranking.com/coin/omwkOTglq+elrond-egld'; D='30.82'; E='  +1.78%  '},
    @{Row=50; B='NEARProtocol'; C='https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'; D='1.296'; E='  +5.04%  '},
    @{Row=51; B='Aave'; C='https://coinranking.com/coin/ixgUfzmLR+aave-aave'; D='53.07'; E='  +2.63%  '}
)

foreach ($u in $updates) {
    $ws.Cells.Item($u.Row, 2).Value = $u.B
    $ws.Cells.Item($u.Row, 3).Value = $u.C

    # Price/Volume text (e.g. "28.320.73", "0.9988") can look numeric to
    # Excels auto-detection; force Text format so COM keeps them as
    # inline strings (matching the source feed formatting), then restore
    # the default "Normal" style so no stray per-cell format lingers.
    $dCell = $ws.Cells.Item($u.Row, 4)
    $dCell.NumberFormat = "@"
    $dCell.Value = $u.D
    $dCell.Style = "Normal"

    $eCell = $ws.Cells.Item($u.Row, 5)
    $eCell.NumberFormat = "@"
    $eCell.Value = $u.E
    $eCell.Style = "Normal"
}

Write-Output "Updated $($updates.Count) rows"